# Relabel the header row with math-mode variable names.
# Written D1, C1, B1 (reverse column order) so the new shared-string
# entries land in the workbook in the order $b$, $A$, $a$.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "`$b`$"
$ws.Range("C1").Value = "`$A`$"
$ws.Range("B1").Value = "`$a`$"

# Select the whole data range and zoom out a bit.
$ws.Range("A1:G6").Select()

$window = $excel.ActiveWindow
$window.Zoom = 225
